$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (Date) values between row 2 and row 5
$ws.Range("D2").Value = 44980
$ws.Range("D5").Value = 44981

# Swap the "Volumen" values between row 2 and row 5
$ws.Range("M2").Value = 50
$ws.Range("M5").Value = 30
